$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 19: new issue noted as Resolved, highlighted like the "Check interval" row (red fill)
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial($xlPasteFormats)
$ws.Range("A19").Value = 42608

$ws.Range("B3").Copy()
$ws.Range("B19").PasteSpecial($xlPasteFormats)
$ws.Range("B19").Value = "Notes stay sharp or flat at calibration"

$ws.Range("C19").Value = "Resolved"

# Row 20: new issue, highlighted like the "Add new audio and loops" row (yellow fill), with a note
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)
$ws.Range("A20").Value = 42608

$ws.Range("B4").Copy()
$ws.Range("B20").PasteSpecial($xlPasteFormats)
$ws.Range("B20").Value = "New chord - melody + bass chord + up one third"

$ws.Range("C20").Value = "Resolved"
$ws.Range("D20").Value = "not sure if it's right"

$ws.Range("B20").Select()
